$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 15th row of audit data (row 15: Catégorie / Problème / Explication / Bonne pratique)
# Shared-string table must grow in this exact order: explication, bonne pratique, problème.
$ws.Range("C15").Value = "Répéter les mêmes mots clés est une pratique pénalisé par Google"
$ws.Range("D15").Value = "Modifier les mots clés répétés par des mots clés pertinants"
$ws.Range("B15").Value = "balise meta keywords à modifier + supprimer les balises méta keywords ittérées dans la page html"
$ws.Range("A15").Value = "SEO"

# Copy the style used by the previous data rows (A14:D14) onto the new row
$ws.Range("A14:D14").Copy() | Out-Null
$ws.Range("A15:D15").PasteSpecial(-4122) | Out-Null

# Update the view: scroll back to column A and move the selection to B15
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B15").Select() | Out-Null
